$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 123, pushing the existing rows 123-131 down to 124-132.
$ws.Rows("123:123").Insert()

# Populate the newly inserted row 123 with a new weekly price record.
# Columns A,B,C,E-L,Q,R,T mirror the constant "metadata" values used across
# this product's rows; D,M,N,O,P,S carry the new observation's data.
$ws.Range("A123").Value = 10
$ws.Range("B123").Value = "Vega Modelo de Temuco"
$ws.Range("C123").Value = "La Araucanía"
$ws.Range("D123").Value = 45265
$ws.Range("E123").Value = 9
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100108
$ws.Range("H123").Value = "Tropicales y subtropicales"
$ws.Range("I123").Value = 100108007
$ws.Range("J123").Value = "Coco"
$ws.Range("K123").Value = "Sin especificar"
$ws.Range("L123").Value = "Primera"
$ws.Range("M123").Value = 50
$ws.Range("N123").Value = 30000
$ws.Range("O123").Value = 30000
$ws.Range("P123").Value = 30000
$ws.Range("Q123").Value = "$/malla 20 unidades"
$ws.Range("R123").Value = "Perú"
$ws.Range("S123").Value = 1500
$ws.Range("T123").Value = 20
